# Small changes in AVL classes
# Swap the SFORZA / TORENBEEK_1982 rows (label + value together) in the
# Xcg/Ycg "ESTIMATION METHOD COMPARISON" tables on the FUSELAGE and WING
# sheets.

$wb = $excel.ActiveWorkbook

# --- FUSELAGE sheet: rows 23 (SFORZA) and 24 (TORENBEEK_1982) swap places ---
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")

$label23 = $wsFuselage.Range("A23").Value2
$value23 = $wsFuselage.Range("C23").Value2
$label24 = $wsFuselage.Range("A24").Value2
$value24 = $wsFuselage.Range("C24").Value2

$wsFuselage.Range("A23").Value = $label24
$wsFuselage.Range("C23").Value = $value24
$wsFuselage.Range("A24").Value = $label23
$wsFuselage.Range("C24").Value = $value23

# --- WING sheet: rows 23/24 (Xcg) and rows 27/28 (Ycg) swap places ---
$wsWing = $wb.Worksheets.Item("WING")

$label23 = $wsWing.Range("A23").Value2
$value23 = $wsWing.Range("C23").Value2
$label24 = $wsWing.Range("A24").Value2
$value24 = $wsWing.Range("C24").Value2

$wsWing.Range("A23").Value = $label24
$wsWing.Range("C23").Value = $value24
$wsWing.Range("A24").Value = $label23
$wsWing.Range("C24").Value = $value23

$label27 = $wsWing.Range("A27").Value2
$value27 = $wsWing.Range("C27").Value2
$label28 = $wsWing.Range("A28").Value2
$value28 = $wsWing.Range("C28").Value2

$wsWing.Range("A27").Value = $label28
$wsWing.Range("C27").Value = $value28
$wsWing.Range("A28").Value = $label27
$wsWing.Range("C28").Value = $value27
